$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the confidential disclosure date text in A7
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-13 for illustrative purposes only and are subject to change."

# Update the weight/percent-change values
$ws.Range("D2").Value = 0.8482553613136102
$ws.Range("E2").Value = -0.005256570713391695

$ws.Range("D3").Value = 0.1517446386863897
$ws.Range("E3").Value = 0.000932835820895539

$ws.Range("D4").Value = 0.9999999999999999
$ws.Range("E4").Value = -0.00431736145516326
